# Major accuracy check update:
#  - The PolyA isolation kit string "NEBNextPoly(A)E7490" was actually the
#    "L" SKU, so every row in column G (polyAIsolationProtocol) should read
#    "NEBNextPoly(A)E7490L" instead of the per-row catalog numbers that had
#    crept in for rows 22-27.
#  - The selection cursor moves from column I to column G.
#  - Columns G/H/I get custom widths to better display the corrected text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$correctKit = "NEBNextPoly(A)E7490L"

# Fix the drifted catalog numbers in column G (rows 2-27) so every row uses
# the single, accurate kit string.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 7).Value = $correctKit
}

# Resize columns G, I to their reviewed widths (H keeps the sheet default).
$ws.Columns.Item(7).ColumnWidth = 18.4167
$ws.Columns.Item(9).ColumnWidth = 16.9167

# Move the active selection from I2:I27 to G2:G27.
$ws.Range("G2:G27").Select()
